# Kearsley_GSP_group_only.xlsx - "line&trafo" sheet edits
# Fill in the previously-blank "Length (km)" column (H) with 0 for the
# rows that didn't have a recorded length (transformers), and leave the
# active selection on the last-edited cell (O15), matching the author's
# manual data-entry pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("line&trafo")
$ws.Activate()

# Rows 10-17 and 50-81 previously had an empty column-H cell; set them to 0.
$rowsToFill = @(10,11,12,13,14,15,16,17) + (50..81)

foreach ($r in $rowsToFill) {
    $ws.Cells.Item($r, 8).Value = 0
}

# Leave the selection where the user's last edit was (O15), matching the
# sheetView's <selection> change in the saved file.
$ws.Range("O15").Select()
